$d = $word.ActiveDocument

# Locate the paragraph that begins "Vous allez participer..." (the
# Globe at Night intro paragraph) so the edit is robust to paragraph
# numbering.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Vous allez*") {
        $target = $p
        break
    }
}

$newText = "Vous allez participer à une campagne mondiale d’observation pour détecter les plus faibles étoiles visibles afin de mesurer la pollution lumineuse sur un site donné. Partout dans le monde, en localisant et en observant la Constellation du Lion dans le ciel nocturne et en la comparant aux cartes stellaires, les participants, apprendront comment l’éclairage, dans leur environnement local, influence la pollution lumineuse. Vos contributions à la base de données en ligne permettront de mesurer la qualité du ciel nocturne."

$r = $target.Range
[void]$r.MoveEnd(1, -1)    # exclude the paragraph mark
$r.Delete()                # remove all the old runs
$r.InsertAfter($newText)   # insert the new text as a single, unformatted run
